# Update the Thrones Pool answer key worksheet for episode 5 ("The Bells")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Characters who died in episode 5 -> flip their answer from "Lives" to "Dies"
$ws.Range("E5").Value = "Dies"   # Sandor Clegane
$ws.Range("C7").Value = "Dies"   # Euron Greyjoy
$ws.Range("E10").Value = "Dies"  # Jaime Lannister
$ws.Range("E13").Value = "Dies" # Gregor Clegane
$ws.Range("E18").Value = "Dies" # Cercei Lannister
$ws.Range("E22").Value = "Dies" # Varys
$ws.Range("E23").Value = "Dies" # Qyburn

# Bonus question answers revealed by episode 5 events
$ws.Range("E28").Value = "Sandor Clegane"      # Winner of CLEGANEBOWL?
$ws.Range("E30").Value = "Daenerys Targaryen"  # Who kills Cercei?

# Update the active view/selection state to match the author's saved session
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("J31").Select()
